$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H86").Value = 1530.2
$ws.Range("I86").Value = 1074.5
$ws.Range("J86").Value = 1834
$ws.Range("K86").Value = 1074.5
$ws.Range("L86").Value = 1834
$ws.Range("M86").Value = 48.5
$ws.Range("N86").Value = -4080

$ws.Range("H89").Value = 1530.2
$ws.Range("I89").Value = 1074.5
$ws.Range("J89").Value = 1834
$ws.Range("K89").Value = 5372.5
$ws.Range("L89").Value = 9170
$ws.Range("M89").Value = 243.5
$ws.Range("N89").Value = -20402

$ws.Range("H137").Value = 2300.6365
$ws.Range("J137").Value = 2704.6
$ws.Range("L137").Value = 8113.799999999999
$ws.Range("N137").Value = -13213.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 56000000
$ws.Range("J12").Value = 56000000
$ws.Range("L12").Value = 56000000
$ws.Range("N12").Value = -56000346

$ws.Range("H61").Value = 2878.1365
$ws.Range("J61").Value = 6214.2856
$ws.Range("L61").Value = 6214.2856
$ws.Range("N61").Value = -6638.2856

$ws.Range("H74").Value = 2116.75
$ws.Range("I74").Value = 1883.4117
$ws.Range("J74").Value = 3439
$ws.Range("K74").Value = 1883.4117
$ws.Range("L74").Value = 3439
$ws.Range("M74").Value = -1009.4117
$ws.Range("N74").Value = -5187

$ws.Range("H77").Value = 2116.75
$ws.Range("I77").Value = 1883.4117
$ws.Range("J77").Value = 3439
$ws.Range("K77").Value = 9417.058500000001
$ws.Range("L77").Value = 17195
$ws.Range("M77").Value = -5049.058500000001
$ws.Range("N77").Value = -25931

$ws.Range("H122").Value = 2740.625
$ws.Range("I122").Value = 1582.2
$ws.Range("J122").Value = 4671.3335
$ws.Range("K122").Value = 4746.6
$ws.Range("L122").Value = 14014.0005
$ws.Range("M122").Value = -2296.6
$ws.Range("N122").Value = -18914.0005

$ws.Range("H132").Value = 1388.4783
$ws.Range("I132").Value = 1354.0476
$ws.Range("J132").Value = 1750
$ws.Range("K132").Value = 4062.142800000001
$ws.Range("L132").Value = 5250
$ws.Range("M132").Value = -1532.142800000001
$ws.Range("N132").Value = -10310

$ws.Range("H136").Value = 2878.1365
$ws.Range("J136").Value = 6214.2856
$ws.Range("L136").Value = 18642.8568
$ws.Range("N136").Value = -23742.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H86").Value = 4087.3076
$ws.Range("I86").Value = 2451.3157
$ws.Range("K86").Value = 2451.3157
$ws.Range("M86").Value = -1328.3157

$ws.Range("H89").Value = 4087.3076
$ws.Range("I89").Value = 2451.3157
$ws.Range("K89").Value = 12256.5785
$ws.Range("M89").Value = -6640.5785

$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 114
$ws.Range("I11").Value = 114
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 114
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 26
$ws.Range("N11").ClearContents()

$ws.Range("H31").Value = 4756.5454
$ws.Range("I31").Value = 2831.7144
$ws.Range("J31").Value = 8125
$ws.Range("K31").Value = 2831.7144
$ws.Range("L31").Value = 8125
$ws.Range("M31").Value = -2536.7144
$ws.Range("N31").Value = -8715

$ws.Range("H34").Value = 4756.5454
$ws.Range("I34").Value = 2831.7144
$ws.Range("J34").Value = 8125
$ws.Range("K34").Value = 2831.7144
$ws.Range("L34").Value = 8125
$ws.Range("M34").Value = -2629.7144
$ws.Range("N34").Value = -8529

$ws.Range("H107").Value = 258.33334
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 258.33334
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 258.33334
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4098.33334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 430.41666
$ws.Range("I6").Value = 95.625
$ws.Range("J6").Value = 1100
$ws.Range("K6").Value = 286.875
$ws.Range("L6").Value = 3300
$ws.Range("M6").Value = -173.875
$ws.Range("N6").Value = -3526

$ws.Range("H131").Value = 2161.6
$ws.Range("I131").Value = 1187.8334
$ws.Range("K131").Value = 3563.5002
$ws.Range("M131").Value = 1476.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 256.94116
$ws.Range("I2").Value = 179.6
$ws.Range("J2").Value = 367.42856
$ws.Range("K2").Value = 179.6
$ws.Range("L2").Value = 367.42856
$ws.Range("M2").Value = -66.59999999999999
$ws.Range("N2").Value = -593.4285600000001

$ws.Range("H26").Value = 28740.334
$ws.Range("J26").Value = 28740.334
$ws.Range("L26").Value = 28740.334
$ws.Range("N26").Value = -29300.334

$ws.Range("H50").Value = 28740.334
$ws.Range("J50").Value = 28740.334
$ws.Range("L50").Value = 28740.334
$ws.Range("N50").Value = -29736.334

$ws.Range("H92").Value = 9749.25
$ws.Range("J92").Value = 9749.25
$ws.Range("L92").Value = 9749.25
$ws.Range("N92").Value = -13493.25

$ws.Range("H122").Value = 3519.1667
$ws.Range("I122").Value = 3266.6667
$ws.Range("J122").Value = 3771.6667
$ws.Range("K122").Value = 9800.000100000001
$ws.Range("L122").Value = 11315.0001
$ws.Range("M122").Value = -7350.000100000001
$ws.Range("N122").Value = -16215.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3671.6667
$ws.Range("I7").Value = 2507.5
$ws.Range("K7").Value = 2507.5
$ws.Range("M7").Value = -2395.5

$ws.Range("H40").Value = 6299.7856
$ws.Range("I40").Value = 6375.3335
$ws.Range("J40").Value = 5846.5
$ws.Range("K40").Value = 6375.3335
$ws.Range("L40").Value = 5846.5
$ws.Range("M40").Value = -6239.3335
$ws.Range("N40").Value = -6118.5

$ws.Range("H55").Value = 1268.7693
$ws.Range("I55").Value = 524.1667
$ws.Range("K55").Value = 524.1667
$ws.Range("M55").Value = -351.1667

$ws.Range("H56").Value = 13620.833
$ws.Range("I56").Value = 11345
$ws.Range("K56").Value = 11345
$ws.Range("M56").Value = -10654

$ws.Range("H100").Value = 4115.231
$ws.Range("I100").Value = 749.8333
$ws.Range("K100").Value = 749.8333
$ws.Range("M100").Value = -208.8333

$ws.Range("H122").Value = 3333.2
$ws.Range("I122").Value = 3370.2222
$ws.Range("K122").Value = 10110.6666
$ws.Range("M122").Value = -7660.6666

$ws.Range("H126").Value = 3671.6667
$ws.Range("I126").Value = 2507.5
$ws.Range("K126").Value = 7522.5
$ws.Range("M126").Value = -5052.5

$ws.Range("H132").Value = 5602.8696
$ws.Range("I132").Value = 5326.952
$ws.Range("J132").Value = 8500
$ws.Range("K132").Value = 15980.856
$ws.Range("L132").Value = 25500
$ws.Range("M132").Value = -13450.856
$ws.Range("N132").Value = -30560

$ws.Range("H136").Value = 1282
$ws.Range("I136").Value = 1108
$ws.Range("K136").Value = 3324
$ws.Range("M136").Value = -774

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H48").Value = 10065
$ws.Range("J48").Value = 10065
$ws.Range("L48").Value = 10065
$ws.Range("N48").Value = -11203
